# ------------------------------------------------------------------------
# Heroes of Pymoli report update
#   1. Title line: "Week 1 Homework" -> Pandas Homework "Heroes of Pymoli"
#      and date "11/8/2019" -> "12/5/2019", with the date split across
#      several runs and the "_GoBack" bookmark relocated here.
#   2. Merge the "1967.64" / " with an average..." runs into one run.
#   3. Remove the old "_GoBack" bookmark location and merge the two runs
#      of the "Based on the above data..." paragraph into one run.
# ------------------------------------------------------------------------

$d = $word.ActiveDocument
$lq = [string][char]0x201C
$rq = [string][char]0x201D

# =========================================================================
# Part 0: the document starts with a "_GoBack" bookmark sitting between the
# two runs of the "Based on the above data..." paragraph; it needs to move
# up to the title paragraph. Remove it now so Part 1 is free to add a
# bookmark of the same name in its new home.
# =========================================================================

$d.Bookmarks("_GoBack").Delete()

# =========================================================================
# Part 1: title paragraph
# =========================================================================

$titlePara = $d.Paragraphs.Item(1)
$titleText = $titlePara.Range.Text
$paraStart = $titlePara.Range.Start

$oldTitlePhrase = "Week 1 Homework"
$oldDatePhrase  = "  11/8/2019"

$titlePhraseStart = $paraStart + $titleText.IndexOf($oldTitlePhrase)
$titlePhraseEnd   = $titlePhraseStart + $oldTitlePhrase.Length

# Guard the boundaries on either side of the run holding "Week 1 Homework"
# so the merge-on-edit pass does not cascade into the neighboring
# "David Winton " / "Written Report " runs, which must stay untouched.
$d.Bookmarks.Add("zzGuardBefore", $d.Range($titlePhraseStart, $titlePhraseStart))
$d.Bookmarks.Add("zzGuardAfter", $d.Range($titlePhraseEnd, $titlePhraseEnd))

# Replace the title wording.
$newTitlePhrase = "Pandas Homework " + $lq + "Heroes of Pymoli" + $rq
$d.Range($titlePhraseStart, $titlePhraseEnd).Text = $newTitlePhrase

# Re-locate the date text (position shifted after the edit above) and
# replace it with the new date.
$titleText = $titlePara.Range.Text
$datePhraseStart = $paraStart + $titleText.IndexOf($oldDatePhrase)
$datePhraseEnd   = $datePhraseStart + $oldDatePhrase.Length
$newDatePhrase = "  12/5/2019"
$d.Range($datePhraseStart, $datePhraseEnd).Text = $newDatePhrase

# Boundaries are now safely fixed; drop the guard bookmarks.
$d.Bookmarks("zzGuardBefore").Delete()
$d.Bookmarks("zzGuardAfter").Delete()

# Re-derive the position right after the new title text (this is where the
# relocated "_GoBack" bookmark belongs) and the start of the new date text.
$titleText = $titlePara.Range.Text
$afterTitlePos = $paraStart + $titleText.IndexOf($newDatePhrase)

# Split "  12/5/2019" into "  1" | "2" | "/" | "5" | "/2019" by touching
# (adding then removing) bookmarks at each internal boundary -- the runs
# stay split even after the marker bookmark is deleted.
$splitOffsets = @(3, 4, 5, 6)
$i = 0
foreach ($off in $splitOffsets) {
    $pos = $afterTitlePos + $off
    $markName = "zzSplit" + $i
    $d.Bookmarks.Add($markName, $d.Range($pos, $pos))
    $i = $i + 1
}
$i = 0
foreach ($off in $splitOffsets) {
    $markName = "zzSplit" + $i
    $d.Bookmarks($markName).Delete()
    $i = $i + 1
}

# Finally, place the real "_GoBack" bookmark right after the title text and
# before the (now split) date runs, collapsed (start == end).
$d.Bookmarks.Add("_GoBack", $d.Range($afterTitlePos, $afterTitlePos))

# =========================================================================
# Part 2: merge "1967.64" and " with an average purchase..." into one run
# =========================================================================

$findRange1 = $d.Content
$found = $findRange1.Find.Execute("1967.64 with an average", $true, $false, $false,
                                   $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $hitStart = $findRange1.Start
    $mergePos = $hitStart + "1967.64".Length
    $r = $d.Range($mergePos, $mergePos)
    $r.InsertAfter("x")
    $d.Range($mergePos, $mergePos + 1).Text = ""
}

# =========================================================================
# Part 3: merge the two runs of the "Based on the above data..." paragraph
# (formerly split by the now-relocated "_GoBack" bookmark) into one run.
# =========================================================================

$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute("be researched for", $true, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $pos2 = $findRange2.Start
    $r2 = $d.Range($pos2, $pos2)
    $r2.InsertAfter("x")
    $d.Range($pos2, $pos2 + 1).Text = ""
}
